# Auto-generated Excel COM-interop script
# Applies scheduled-runner market-price/profit updates to the Mateus_Profits workbook.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4399.4
$ws.Range("J51").Value = 3666.6667
$ws.Range("L51").Value = 3666.6667
$ws.Range("N51").Value = -4634.6667
$ws.Range("H75").Value = 72999.664
$ws.Range("J75").Value = 72999.664
$ws.Range("L75").Value = 72999.664
$ws.Range("N75").Value = -74871.664
$ws.Range("H78").Value = 72999.664
$ws.Range("J78").Value = 72999.664
$ws.Range("L78").Value = 218998.992
$ws.Range("N78").Value = -228358.992
$ws.Range("H98").Value = 2614.6428
$ws.Range("I98").Value = 2614.6428
$ws.Range("K98").Value = 2614.6428
$ws.Range("M98").Value = -1116.6428
$ws.Range("H111").Value = 766.44446
$ws.Range("I111").Value = 799.3333
$ws.Range("J111").Value = 750
$ws.Range("K111").Value = 2397.9999
$ws.Range("L111").Value = 2250
$ws.Range("M111").Value = 669.0001000000002
$ws.Range("N111").Value = -8384
$ws.Range("H113").Value = 3390
$ws.Range("I113").Value = 3390
$ws.Range("K113").Value = 3390
$ws.Range("M113").Value = -136
$ws.Range("H116").Value = 4533.3335
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("H122").Value = 2614.6428
$ws.Range("I122").Value = 2614.6428
$ws.Range("K122").Value = 7843.928400000001
$ws.Range("M122").Value = -5393.928400000001
$ws.Range("H132").Value = 995.1667
$ws.Range("I132").Value = 1068.5
$ws.Range("K132").Value = 3205.5
$ws.Range("M132").Value = -675.5
$ws.Range("H135").Value = 2609.375
$ws.Range("I135").Value = 2440.9092
$ws.Range("K135").Value = 21968.1828
$ws.Range("M135").Value = -19433.1828
$ws.Range("H137").Value = 1594.4324
$ws.Range("I137").Value = 1488.6538
$ws.Range("J137").Value = 1844.4546
$ws.Range("K137").Value = 4465.9614
$ws.Range("L137").Value = 5533.3638
$ws.Range("M137").Value = -1915.9614
$ws.Range("N137").Value = -10633.3638
$ws.Range("H138").Value = 2388.2144
$ws.Range("I138").Value = 1642.7587
$ws.Range("J138").Value = 3188.889
$ws.Range("K138").Value = 4928.2761
$ws.Range("L138").Value = 9566.667000000001
$ws.Range("M138").Value = 211.7239
$ws.Range("N138").Value = -19846.667
$ws.Range("H141").Value = 2590.75
$ws.Range("I141").Value = 2371.7273
$ws.Range("K141").Value = 7115.1819
$ws.Range("M141").Value = -1935.1819

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 1488.25
$ws.Range("J26").Value = 1466.3334
$ws.Range("L26").Value = 1466.3334
$ws.Range("N26").Value = -2126.3334
$ws.Range("H27").Value = 966.3333
$ws.Range("J27").Value = 966.3333
$ws.Range("L27").Value = 966.3333
$ws.Range("N27").Value = -1334.3333
$ws.Range("H32").Value = 3665.463
$ws.Range("I32").Value = 3842.392
$ws.Range("K32").Value = 3842.392
$ws.Range("M32").Value = -3555.392
$ws.Range("H74").Value = 1817.174
$ws.Range("I74").Value = 1743.4857
$ws.Range("J74").Value = 2051.6365
$ws.Range("K74").Value = 1743.4857
$ws.Range("L74").Value = 2051.6365
$ws.Range("M74").Value = -869.4857
$ws.Range("N74").Value = -3799.6365
$ws.Range("H77").Value = 1817.174
$ws.Range("I77").Value = 1743.4857
$ws.Range("J77").Value = 2051.6365
$ws.Range("K77").Value = 8717.4285
$ws.Range("L77").Value = 10258.1825
$ws.Range("M77").Value = -4349.4285
$ws.Range("N77").Value = -18994.1825
$ws.Range("H132").Value = 7248.8335
$ws.Range("I132").Value = 5743.636
$ws.Range("K132").Value = 17230.908
$ws.Range("M132").Value = -14700.908

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 34999.5
$ws.Range("J103").Value = 34999.5
$ws.Range("L103").Value = 34999.5
$ws.Range("N103").Value = -37343.5
$ws.Range("H134").Value = 2167.3704
$ws.Range("I134").Value = 2215.0613
$ws.Range("K134").Value = 6645.1839
$ws.Range("M134").Value = -4110.1839
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 2998.3333
$ws.Range("I6").Value = 2998.3333
$ws.Range("K6").Value = 2998.3333
$ws.Range("M6").Value = -2885.3333
$ws.Range("H7").Value = 292.92307
$ws.Range("I7").Value = 330.27274
$ws.Range("K7").Value = 330.27274
$ws.Range("M7").Value = -217.27274
$ws.Range("H31").Value = 7400.2144
$ws.Range("I31").Value = 5688.8887
$ws.Range("J31").Value = 10480.6
$ws.Range("K31").Value = 5688.8887
$ws.Range("L31").Value = 10480.6
$ws.Range("M31").Value = -5393.8887
$ws.Range("N31").Value = -11070.6
$ws.Range("H34").Value = 7400.2144
$ws.Range("I34").Value = 5688.8887
$ws.Range("J34").Value = 10480.6
$ws.Range("K34").Value = 5688.8887
$ws.Range("L34").Value = 10480.6
$ws.Range("M34").Value = -5486.8887
$ws.Range("N34").Value = -10884.6
$ws.Range("H96").Value = 52437
$ws.Range("J96").Value = 52437
$ws.Range("L96").Value = 52437
$ws.Range("N96").Value = -57929
$ws.Range("H99").Value = 5806.25
$ws.Range("J99").Value = 9450
$ws.Range("L99").Value = 9450
$ws.Range("N99").Value = -12446
$ws.Range("H126").Value = 5806.25
$ws.Range("J126").Value = 9450
$ws.Range("L126").Value = 28350
$ws.Range("N126").Value = -33290
$ws.Range("H141").Value = 39900
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 431.8
$ws.Range("J122").Value = 471.25
$ws.Range("L122").Value = 4241.25
$ws.Range("N122").Value = -9141.25
$ws.Range("H132").Value = 1479.1538
$ws.Range("J132").Value = 1697.375
$ws.Range("L132").Value = 15276.375
$ws.Range("N132").Value = -20336.375

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2718.76
$ws.Range("I126").Value = 1824.7693
$ws.Range("K126").Value = 5474.3079
$ws.Range("M126").Value = -3004.3079
$ws.Range("H132").Value = 2472.383
$ws.Range("I132").Value = 2260.9412
$ws.Range("K132").Value = 6782.823600000001
$ws.Range("M132").Value = -4252.823600000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H76").Value = 17696.75
$ws.Range("J76").Value = 17696.75
$ws.Range("L76").Value = 17696.75
$ws.Range("N76").Value = -18372.75
$ws.Range("H79").Value = 17696.75
$ws.Range("J79").Value = 17696.75
$ws.Range("L79").Value = 17696.75
$ws.Range("N79").Value = -20036.75
$ws.Range("H93").Value = 21466.334
$ws.Range("I93").Value = 1165.8334
$ws.Range("K93").Value = 1165.8334
$ws.Range("M93").Value = 82.16660000000002
$ws.Range("H136").Value = 9232.799999999999
$ws.Range("I136").Value = 9139
$ws.Range("J136").Value = 9490.75
$ws.Range("K136").Value = 27417
$ws.Range("L136").Value = 28472.25
$ws.Range("M136").Value = -24867
$ws.Range("N136").Value = -33572.25

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6399.3335
$ws.Range("J62").Value = 6399.3335
$ws.Range("L62").Value = 6399.3335
$ws.Range("N62").Value = -7647.3335
$ws.Range("H63").Value = 50150.668
$ws.Range("J63").Value = 50150.668
$ws.Range("L63").Value = 50150.668
$ws.Range("N63").Value = -51398.668
$ws.Range("H65").Value = 6399.3335
$ws.Range("J65").Value = 6399.3335
$ws.Range("L65").Value = 31996.6675
$ws.Range("N65").Value = -38236.6675
$ws.Range("H66").Value = 50150.668
$ws.Range("J66").Value = 50150.668
$ws.Range("L66").Value = 150452.004
$ws.Range("N66").Value = -156692.004
$ws.Range("H81").Value = 2550.3333
$ws.Range("I81").Value = 1516
$ws.Range("K81").Value = 3032
$ws.Range("M81").Value = -1971
$ws.Range("H82").Value = 65000
$ws.Range("J82").Value = 65000
$ws.Range("L82").Value = 65000
$ws.Range("N82").Value = -65766
$ws.Range("H84").Value = 2550.3333
$ws.Range("I84").Value = 1516
$ws.Range("K84").Value = 15160
$ws.Range("M84").Value = -9856
$ws.Range("H85").Value = 65000
$ws.Range("J85").Value = 65000
$ws.Range("L85").Value = 65000
$ws.Range("N85").Value = -67652
$ws.Range("H126").Value = 2209.25
$ws.Range("I126").Value = 2534.5386
$ws.Range("K126").Value = 7603.6158
$ws.Range("M126").Value = -5133.6158
$ws.Range("H136").Value = 5342.6895
$ws.Range("I136").Value = 3443.4707
$ws.Range("K136").Value = 10330.4121
$ws.Range("M136").Value = -7780.4121
